$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1063
$ws.Range("F3").Value = 671
$ws.Range("F4").Value = 1482
$ws.Range("F5").Value = 3234
$ws.Range("F7").Value = 648
$ws.Range("F8").Value = 2207
$ws.Range("F9").Value = 474
$ws.Range("F10").Value = 405
$ws.Range("F12").Value = 125
$ws.Range("F13").Value = 305
$ws.Range("F14").Value = 1063
$ws.Range("F15").Value = 426
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 195
$ws.Range("F19").Value = 4402
$ws.Range("F20").Value = 1285
$ws.Range("F21").Value = 3361
$ws.Range("F23").Value = 58
$ws.Range("F24").Value = 158
$ws.Range("F25").Value = 3277
$ws.Range("F26").Value = 4884
$ws.Range("F29").Value = 539
$ws.Range("F30").Value = 3164
$ws.Range("F31").Value = 334
$ws.Range("F33").Value = 129
$ws.Range("F35").Value = 870
$ws.Range("F36").Value = 1144
$ws.Range("F37").Value = 1387
$ws.Range("F39").Value = 1314
$ws.Range("F40").Value = 838
$ws.Range("F42").Value = 786
$ws.Range("F44").Value = 50
$ws.Range("F45").Value = 278
$ws.Range("F46").Value = 56
$ws.Range("F47").Value = 139
$ws.Range("F49").Value = 3707

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 993

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2082

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2082
$ws.Range("F3").Value = 671
$ws.Range("F4").Value = 1482
$ws.Range("F5").Value = 3234
$ws.Range("F7").Value = 648
$ws.Range("F9").Value = 2207
$ws.Range("F10").Value = 474
$ws.Range("F11").Value = 405
$ws.Range("F13").Value = 993
$ws.Range("F14").Value = 125
$ws.Range("F15").Value = 305
$ws.Range("F16").Value = 1063
$ws.Range("F17").Value = 426
$ws.Range("F18").Value = 195
$ws.Range("F19").Value = 4402
$ws.Range("F21").Value = 1285
$ws.Range("F23").Value = 3361
$ws.Range("F24").Value = 3277
$ws.Range("F25").Value = 4884
$ws.Range("F28").Value = 3164
$ws.Range("F29").Value = 334
$ws.Range("F31").Value = 129
$ws.Range("F33").Value = 870
$ws.Range("F34").Value = 1144
$ws.Range("F35").Value = 1387
$ws.Range("F37").Value = 1314
$ws.Range("F39").Value = 838
$ws.Range("F42").Value = 50
$ws.Range("F44").Value = 278
$ws.Range("F46").Value = 56
$ws.Range("F47").Value = 139
$ws.Range("F49").Value = 3707

Write-Output "Update complete"
